# "Diseño Elastic Flex-wrap max-width"
#
# The tracked changes are mostly cosmetic proofing artifacts left behind by
# Word's spelling/grammar checker (<w:proofErr .../> run-splits) plus one
# real text fix ("varian" -> "varían") and the "_GoBack" bookmark moving
# from the last paragraph up to the first one. We rebuild each affected
# paragraph's content via Range.InsertXML so the exact run/proofErr layout
# from the target matches.

$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Paragraph 1: "Unidades relativas de medida: " -> the _GoBack bookmark moves here ---
$p1 = $d.Paragraphs(1).Range
$null = $p1.InsertXML(@"
<w:p $wns>
  <w:r w:rsidRPr="008A33BE">
    <w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr>
    <w:t>Unidades relativas de medida:</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
"@)

# --- Paragraph 2: "varian..." -> "varían..." (real fix), split into two runs ---
$p2 = $d.Paragraphs(2).Range
$null = $p2.InsertXML(@"
<w:p $wns>
  <w:r>
    <w:t>var&#237;an</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> seg&#250;n alguna condici&#243;n</w:t>
  </w:r>
</w:p>
"@)

# --- Paragraph 3: "Porcentaje % : " -> split with gramStart/gramEnd proofErr around "% :" ---
$p3 = $d.Paragraphs(3).Range
$null = $p3.InsertXML(@"
<w:p $wns>
  <w:r>
    <w:t xml:space="preserve">Porcentaje </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:t>% :</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
</w:p>
"@)

# --- Paragraph 4: "Longitud referente al tamaño de los elementos padre." unchanged ---
# --- Paragraph 5: "Em:" unchanged ---

# --- Paragraph 6: wrap "mas" with spellStart/spellEnd proofErr ---
$p6 = $d.Paragraphs(6).Range
$null = $p6.InsertXML(@"
<w:p $wns>
  <w:r>
    <w:t xml:space="preserve"> unidad relativa al tama&#241;o de fuente especificada </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>mas</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> cercana.</w:t>
  </w:r>
</w:p>
"@)

# --- Paragraph 7: "Rem:" unchanged ---

# --- Paragraph 8: wrap "mas" / "html" / "body" with spellStart/spellEnd proofErr ---
$p8 = $d.Paragraphs(8).Range
$null = $p8.InsertXML(@"
<w:p $wns>
  <w:r>
    <w:t xml:space="preserve">Unidad relativa al tama&#241;o de fuente especificada en el ancestro </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>mas</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> lejano (</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>html</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> o </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>body</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>)</w:t>
  </w:r>
  <w:r>
    <w:t>.</w:t>
  </w:r>
</w:p>
"@)

# --- Paragraph 9: "Vw/vh:" -> wrap "Vw" and "vh" with spellStart/spellEnd proofErr ---
$p9 = $d.Paragraphs(9).Range
$null = $p9.InsertXML(@"
<w:p $wns>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>Vw</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>/</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>vh</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>:</w:t>
  </w:r>
</w:p>
"@)

# --- Paragraph 10: wrap "Viewport" with spellStart/spellEnd; bookmark removed (now on paragraph 1) ---
$p10 = $d.Paragraphs(10).Range
$null = $p10.InsertXML(@"
<w:p $wns>
  <w:r>
    <w:t xml:space="preserve">Unidad relativa porcentual con respecto al </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>Viewport</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>.</w:t>
  </w:r>
</w:p>
"@)
